# power bi prototipo evaluacion
# Adds a new "Entregas" sheet and reworks the "Evaluaciones" sheet
# (new columns/rows for the evaluation-vs-exercise breakdown).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluaciones")

# --- Step 1: rename the E1/F1 headers on Evaluaciones --------------------
# (must happen before the Entregas sheet is populated so the shared-string
# table is built up in the same order Excel would produce it)
$ws.Range("E1").Value = "promedioAciertos"
$ws.Range("F1").Value = "promedioFallos"

# --- Step 2: update existing data rows (D/E/F columns) -------------------
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 814
$ws.Range("F2").Value = 53

$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 747
$ws.Range("F3").Value = 96

$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 400
$ws.Range("F4").Value = 13

$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 440
$ws.Range("F5").Value = 84

# --- Step 3: create & populate the new "Entregas" sheet -------------------
# inserted right after "Evaluaciones" (and before the hidden "Preguntas")
$evalSheet = $wb.Worksheets.Item("Evaluaciones")
$entregas = $wb.Worksheets.Add($null, $evalSheet)
$entregas.Name = "Entregas"

$entregas.Range("A1").Value = "Ejercicio"
$entregas.Range("B1").Value = "Correcto"
$entregas.Range("C1").Value = "Incorrecto"
$entregas.Range("D1").Value = "Total"

$entregas.Range("A2").Value = "Ejercicio 1"
$entregas.Range("B2").Value = 82
$entregas.Range("C2").Value = 18
$entregas.Range("D2").Value = 100

$entregas.Range("A3").Value = "Ejercicio 2"
$entregas.Range("B3").Value = 91
$entregas.Range("C3").Value = 9
$entregas.Range("D3").Value = 100

$entregas.Range("A4").Value = "Ejercicio 3"
$entregas.Range("B4").Value = 85
$entregas.Range("C4").Value = 15
$entregas.Range("D4").Value = 100

$entregas.Range("A5").Value = "Ejercicio 4"
$entregas.Range("B5").Value = 92
$entregas.Range("C5").Value = 8
$entregas.Range("D5").Value = 100

$entregas.Range("A6").Value = "Ejercicio 5"
$entregas.Range("B6").Value = 40
$entregas.Range("C6").Value = 60
$entregas.Range("D6").Value = 100

$entregas.Range("D9").Select() | Out-Null

# --- Step 4: new G column + new rows 6/7 on Evaluaciones ------------------
$ws.Range("G1").Value = "Evaluacion"
$ws.Range("G2").Value = "Conocimientos Basicos"
$ws.Range("G3").Value = "Introduccion"
$ws.Range("G4").Value = "Practicas"
$ws.Range("G5").Value = "Conocimientos Avanzados"
$ws.Range("G6").Value = "Relacion con otras ciencias"
$ws.Range("G7").Value = "Conclusiones "

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 492
$ws.Range("F6").Value = 151

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 649
$ws.Range("F7").Value = 100

# B7 carries the underlined-font style already used elsewhere in the sheet
$ws.Range("B7").Font.Underline = $true

# --- Step 5: cosmetic touch-ups to line up with the authored layout ------
$ws.Columns.Item(5).ColumnWidth = 15.6
$ws.Columns.Item(7).ColumnWidth = 20.6
$ws.PageSetup.Orientation = 1

# keep "Evaluaciones" as the selected/active tab (it was before the edit)
$evalSheet.Activate() | Out-Null
$ws.Range("F11").Select() | Out-Null
